# Burndown chart update: bump sprint number, rescale goal velocity (E9) and
# update the recorded "Done"/Actual values (F9:F29) for Sprint 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BURNDOWN CHART")

# --- Update the chart title: "Sprint 4" -> "Sprint 5" ---
$chart = $ws.ChartObjects().Item(1).Chart
$chart.ChartTitle.Text = "Withdrive Sprint 5 Burndown chart"

# --- Goal velocity seed (E9). E10:E29 are formulas driven off this value and
#     recalculate automatically. ---
$ws.Range("E9").Value = 40

# --- Done / Actual values (F9:F29) ---
$ws.Range("F9").Value = 40
$ws.Range("F10").Value = 40
$ws.Range("F11").Value = 40
$ws.Range("F12").Value = 40
$ws.Range("F13").Value = 35
$ws.Range("F14").Value = 35
$ws.Range("F15").Value = 35
$ws.Range("F16").Value = 35
$ws.Range("F17").Value = 35
$ws.Range("F18").Value = 20
$ws.Range("F19").Value = 20
$ws.Range("F20").Value = 24
$ws.Range("F21").Value = 24
$ws.Range("F22").Value = 24
$ws.Range("F23").Value = 24
$ws.Range("F24").Value = 20
$ws.Range("F25").Value = 20
$ws.Range("F26").Value = 20
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0

# --- Final selection on the sheet ---
$ws.Range("P32").Select()
